# Update column G ("K") values in the save_data sheet.
# The commit regenerates the save_data export: column G now reports the
# strikeout count "K" (instead of the previous "Strike#" value) for each
# outing, so each value is recalculated/rewritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 1
    6 = 0
    7 = 1
    8 = 0
    9 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 5
    29 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    39 = 3
    40 = 3
    41 = 0
    42 = 0
    43 = 0
    44 = 1
    45 = 0
    46 = 2
    47 = 2
    48 = 1
    49 = 1
    50 = 1
    51 = 1
    52 = 0
    53 = 1
    54 = 0
    55 = 0
    56 = 2
    57 = 0
    58 = 1
    59 = 0
    60 = 0
    62 = 0
    63 = 2
    64 = 2
    65 = 3
    66 = 1
    67 = 0
    68 = 1
    69 = 1
    70 = 0
    71 = 0
    72 = 1
    73 = 0
    77 = 2
    79 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

Write-Host "Updated $($gValues.Count) cells in column G"
